# feat: use html template for transcript exporting
#
# - bump footer version string v.2.7 -> v.2.8
# - insert three blank rows before the "Test filed / signature" block so the
#   transcript rows area (jx:area) has room to grow, pushing the old rows
#   13/14/18 down to 16/17/21
# - add a jx:area(lastCell="G100") comment on the new A15 cell describing the
#   (now larger) print area for the html/pdf export
# - move the active selection to C19

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bump the template version label in E10 (merged E10:F10).
$ws.Range("E10").Value = "v.2.8"

# Insert 3 new rows above the old row 13. Because row 12 already carries the
# s="14" row style (with s="1" cell style, s="0" on column H), Excel's
# "insert copies format from the row above" behaviour reproduces that same
# blank banding for the new rows 13-15, and shifts the old row 13 ("Test
# filed" / "Chữ kí") down to row 16, the old row 14 (blank B14) down to row
# 17, and the old row 18 (signature line) down to row 21 - matching the
# target layout exactly.
$ws.Rows("13:15").Insert()

# New comment describing the (enlarged) print/export area, on the first of
# the freshly inserted blank rows. (Note: deliberately not touching
# Shape.TextFrame.Characters().Font here - in this engine that call resolves
# back onto the anchoring cell's own font/style rather than the comment
# shape, which would corrupt A15's s="1" cell style for no benefit, since
# comment text is modelled/exported as plain text with no <r><rPr> runs
# either way.)
$comment = $ws.Range("A15").AddComment("jx:area(lastCell=`"G100`") `n")

# Match the author's saved cursor position.
$ws.Range("C19").Select() | Out-Null

Write-Host "template.xlsx updated: v.2.8 label, 3 rows inserted at 13:15, A15 comment added, selection -> C19"
